$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.036566368707438
$ws.Cells.Item(2, 4).Value = 1.044708188474221
$ws.Cells.Item(2, 5).Value = 1.035470784259076
$ws.Cells.Item(2, 6).Value = 1.051816545235598
$ws.Cells.Item(2, 9).Value = 1.039149228070396
$ws.Cells.Item(2, 10).Value = 1.041673986803223
$ws.Cells.Item(2, 11).Value = 1.047478479059697
$ws.Cells.Item(2, 12).Value = 1.038267295479748
$ws.Cells.Item(2, 13).Value = 1.054566998531514
$ws.Cells.Item(2, 14).Value = 1.005712725503983

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.037442711288285
$ws.Cells.Item(3, 4).Value = 1.045414660823857
$ws.Cells.Item(3, 5).Value = 1.03621411744198
$ws.Cells.Item(3, 6).Value = 1.052700979514288
$ws.Cells.Item(3, 9).Value = 1.039357220548455
$ws.Cells.Item(3, 10).Value = 1.042194715197283
$ws.Cells.Item(3, 11).Value = 1.047996517510298
$ws.Cells.Item(3, 12).Value = 1.038820205951724
$ws.Cells.Item(3, 13).Value = 1.055263970268134

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.038010380233052
$ws.Cells.Item(4, 4).Value = 1.045872337216553
$ws.Cells.Item(4, 5).Value = 1.036696005294813
$ws.Cells.Item(4, 6).Value = 1.053274276883259
$ws.Cells.Item(4, 9).Value = 1.039490919280576
$ws.Cells.Item(4, 10).Value = 1.042531628192384
$ws.Cells.Item(4, 11).Value = 1.048331589507334
$ws.Cells.Item(4, 12).Value = 1.039178211139985
$ws.Cells.Item(4, 13).Value = 1.055715327550684

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.038249174445167
$ws.Cells.Item(5, 4).Value = 1.04606487234907
$ws.Cells.Item(5, 5).Value = 1.036898805099212
$ws.Cells.Item(5, 6).Value = 1.053515530940068
$ws.Cells.Item(5, 9).Value = 1.039546913373726
$ws.Cells.Item(5, 10).Value = 1.042673257296444
$ws.Cells.Item(5, 11).Value = 1.048472420551437
$ws.Cells.Item(5, 12).Value = 1.039328771524983
$ws.Cells.Item(5, 13).Value = 1.055905165208876

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.03828927757457
$ws.Cells.Item(6, 4).Value = 1.046097207301116
$ws.Cells.Item(6, 5).Value = 1.036932868587551
$ws.Cells.Item(6, 6).Value = 1.053556052570748
$ws.Cells.Item(6, 9).Value = 1.039556302527078
$ws.Cells.Item(6, 10).Value = 1.042697036876725
$ws.Cells.Item(6, 11).Value = 1.04849606471863
$ws.Cells.Item(6, 12).Value = 1.039354054464098
$ws.Cells.Item(6, 13).Value = 1.055937044848084

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.038013570440392
$ws.Cells.Item(7, 4).Value = 1.04587490937965
$ws.Cells.Item(7, 5).Value = 1.036698714276379
$ws.Cells.Item(7, 6).Value = 1.053277499591577
$ws.Cells.Item(7, 9).Value = 1.039491668313592
$ws.Cells.Item(7, 10).Value = 1.042533520684733
$ws.Cells.Item(7, 11).Value = 1.04833347143076
$ws.Cells.Item(7, 12).Value = 1.039180222720708
$ws.Cells.Item(7, 13).Value = 1.055717863830586

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.036862404221456
$ws.Cells.Item(8, 4).Value = 1.044946830979083
$ws.Cells.Item(8, 5).Value = 1.035721809495695
$ws.Cells.Item(8, 6).Value = 1.052115233902201
$ws.Cells.Item(8, 9).Value = 1.039219703215158
$ws.Cells.Item(8, 10).Value = 1.041849975535468
$ws.Cells.Item(8, 11).Value = 1.047653579469023
$ws.Cells.Item(8, 12).Value = 1.03845410469971
$ws.Cells.Item(8, 13).Value = 1.054802465879297

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.03483868964502
$ws.Cells.Item(9, 4).Value = 1.043315659625393
$ws.Cells.Item(9, 5).Value = 1.034007358122027
$ws.Cells.Item(9, 6).Value = 1.050074981385656
$ws.Cells.Item(9, 9).Value = 1.038733707565934
$ws.Cells.Item(9, 10).Value = 1.040645281070026
$ws.Cells.Item(9, 11).Value = 1.04645455601858
$ws.Cells.Item(9, 12).Value = 1.037176449884883
$ws.Cells.Item(9, 13).Value = 1.053192318171863

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.033492851511073
$ws.Cells.Item(10, 4).Value = 1.042231149721478
$ws.Cells.Item(10, 5).Value = 1.032869183910538
$ws.Cells.Item(10, 6).Value = 1.048720168446303
$ws.Cells.Item(10, 9).Value = 1.038405206446973
$ws.Cells.Item(10, 10).Value = 1.039842087741853
$ws.Cells.Item(10, 11).Value = 1.045654631868969
$ws.Cells.Item(10, 12).Value = 1.0363260059432
$ws.Cells.Item(10, 13).Value = 1.052120926600119

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.032910890642273
$ws.Cells.Item(11, 4).Value = 1.041762262701224
$ws.Cells.Item(11, 5).Value = 1.032377500064988
$ws.Cells.Item(11, 6).Value = 1.048134811959934
$ws.Cells.Item(11, 9).Value = 1.03826190193079
$ws.Cells.Item(11, 10).Value = 1.039494296016713
$ws.Cells.Item(11, 11).Value = 1.045308135218164
$ws.Cells.Item(11, 12).Value = 1.035958084985861
$ws.Cells.Item(11, 13).Value = 1.051657505102746

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.032694845474547
$ws.Cells.Item(12, 4).Value = 1.041588205984267
$ws.Cells.Item(12, 5).Value = 1.032195041750406
$ws.Cells.Item(12, 6).Value = 1.047917579565777
$ws.Cells.Item(12, 9).Value = 1.038208513453191
$ws.Cells.Item(12, 10).Value = 1.039365111321366
$ws.Cells.Item(12, 11).Value = 1.045179413532065
$ws.Cells.Item(12, 12).Value = 1.035821473169178
$ws.Cells.Item(12, 13).Value = 1.051485446083825

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.032741182397373
$ws.Cells.Item(13, 4).Value = 1.041625536790355
$ws.Cells.Item(13, 5).Value = 1.032234171733469
$ws.Cells.Item(13, 6).Value = 1.047964167795712
$ws.Cells.Item(13, 9).Value = 1.038219972638293
$ws.Cells.Item(13, 10).Value = 1.039392821829035
$ws.Cells.Item(13, 11).Value = 1.04520702553794
$ws.Cells.Item(13, 12).Value = 1.035850774565655
$ws.Cells.Item(13, 13).Value = 1.051522349847598

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.03289302980443
$ws.Cells.Item(14, 4).Value = 1.041747872889449
$ws.Cells.Item(14, 5).Value = 1.032362414424767
$ws.Cells.Item(14, 6).Value = 1.04811685147277
$ws.Cells.Item(14, 9).Value = 1.038257492058226
$ws.Cells.Item(14, 10).Value = 1.039483617549591
$ws.Cells.Item(14, 11).Value = 1.045297495391146
$ws.Cells.Item(14, 12).Value = 1.035946791574648
$ws.Cells.Item(14, 13).Value = 1.05164328108406

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.032986604081536
$ws.Cells.Item(15, 4).Value = 1.041823262665262
$ws.Cells.Item(15, 5).Value = 1.032441452214403
$ws.Cells.Item(15, 6).Value = 1.048210950827427
$ws.Cells.Item(15, 9).Value = 1.038280587988381
$ws.Cells.Item(15, 10).Value = 1.039539559909805
$ws.Cells.Item(15, 11).Value = 1.045353234585327
$ws.Cells.Item(15, 12).Value = 1.036005957541313
$ws.Cells.Item(15, 13).Value = 1.051717800973936

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.03353149122864
$ws.Cells.Item(16, 4).Value = 1.042262283373601
$ws.Cells.Item(16, 5).Value = 1.032901839789164
$ws.Cells.Item(16, 6).Value = 1.048759043885599
$ws.Cells.Item(16, 9).Value = 1.038414694790167
$ws.Cells.Item(16, 10).Value = 1.039865169556935
$ws.Cells.Item(16, 11).Value = 1.045677625200951
$ws.Cells.Item(16, 12).Value = 1.036350430643252
$ws.Cells.Item(16, 13).Value = 1.052151692983397

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.03387349864135
$ws.Cells.Item(17, 4).Value = 1.042537861529025
$ws.Cells.Item(17, 5).Value = 1.033190938681768
$ws.Cells.Item(17, 6).Value = 1.049103193866339
$ws.Cells.Item(17, 9).Value = 1.038498532750568
$ws.Cells.Item(17, 10).Value = 1.04006941573483
$ws.Cells.Item(17, 11).Value = 1.045881074496793
$ws.Cells.Item(17, 12).Value = 1.036566597951121
$ws.Cells.Item(17, 13).Value = 1.052423996413137

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.034073062477745
$ws.Cells.Item(18, 4).Value = 1.042698670278223
$ws.Cells.Item(18, 5).Value = 1.033359676316823
$ws.Cells.Item(18, 6).Value = 1.049304054731695
$ws.Cells.Item(18, 9).Value = 1.038547331563944
$ws.Cells.Item(18, 10).Value = 1.040188548636308
$ws.Cells.Item(18, 11).Value = 1.045999730959444
$ws.Cells.Item(18, 12).Value = 1.036692716136778
$ws.Cells.Item(18, 13).Value = 1.052582874356244

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.034141121554721
$ws.Cells.Item(19, 4).Value = 1.042753513502813
$ws.Cells.Item(19, 5).Value = 1.033417230289169
$ws.Cells.Item(19, 6).Value = 1.049372564074822
$ws.Cells.Item(19, 9).Value = 1.038563953300856
$ws.Cells.Item(19, 10).Value = 1.040229169766333
$ws.Cells.Item(19, 11).Value = 1.046040187661343
$ws.Cells.Item(19, 12).Value = 1.036735724474253
$ws.Cells.Item(19, 13).Value = 1.052637055717569

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.033836796539026
$ws.Cells.Item(20, 4).Value = 1.042508287485565
$ws.Cells.Item(20, 5).Value = 1.033159909598245
$ws.Cells.Item(20, 6).Value = 1.049066256991142
$ws.Cells.Item(20, 9).Value = 1.03848954832304
$ws.Cells.Item(20, 10).Value = 1.040047502089642
$ws.Cells.Item(20, 11).Value = 1.045859247552089
$ws.Cells.Item(20, 12).Value = 1.03654340198215
$ws.Cells.Item(20, 13).Value = 1.052394775860516

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.032848311175735
$ws.Cells.Item(21, 4).Value = 1.041711844940437
$ws.Cells.Item(21, 5).Value = 1.03232464530506
$ws.Cells.Item(21, 6).Value = 1.048071884536188
$ws.Cells.Item(21, 9).Value = 1.038246447902552
$ws.Cells.Item(21, 10).Value = 1.039456880445731
$ws.Cells.Item(21, 11).Value = 1.045270854744543
$ws.Cells.Item(21, 12).Value = 1.03591851555671
$ws.Cells.Item(21, 13).Value = 1.051607667728943

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.032227512038296
$ws.Cells.Item(22, 4).Value = 1.041211720218887
$ws.Cells.Item(22, 5).Value = 1.031800494594425
$ws.Cells.Item(22, 6).Value = 1.047447812598147
$ws.Cells.Item(22, 9).Value = 1.038092682581326
$ws.Cells.Item(22, 10).Value = 1.039085537453165
$ws.Cells.Item(22, 11).Value = 1.044900808773725
$ws.Cells.Item(22, 12).Value = 1.035525916950625
$ws.Cells.Item(22, 13).Value = 1.051113224150612

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.032556542495643
$ws.Cells.Item(23, 4).Value = 1.041476785485234
$ws.Cells.Item(23, 5).Value = 1.032078260199647
$ws.Cells.Item(23, 6).Value = 1.047778537278768
$ws.Cells.Item(23, 9).Value = 1.038174283324341
$ws.Cells.Item(23, 10).Value = 1.039282392594908
$ws.Cells.Item(23, 11).Value = 1.045096986207887
$ws.Cells.Item(23, 12).Value = 1.035734012760481
$ws.Cells.Item(23, 13).Value = 1.051375295590095

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.033853380402412
$ws.Cells.Item(24, 4).Value = 1.042521650510825
$ws.Cells.Item(24, 5).Value = 1.033173929962263
$ws.Cells.Item(24, 6).Value = 1.04908294679218
$ws.Cells.Item(24, 9).Value = 1.038493608315741
$ws.Cells.Item(24, 10).Value = 1.040057403924136
$ws.Cells.Item(24, 11).Value = 1.045869110246313
$ws.Cells.Item(24, 12).Value = 1.036553883145044
$ws.Cells.Item(24, 13).Value = 1.052407979222456

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.035361292181416
$ws.Cells.Item(25, 4).Value = 1.043736846016449
$ws.Cells.Item(25, 5).Value = 1.034449747398686
$ws.Cells.Item(25, 6).Value = 1.050601499381006
$ws.Cells.Item(25, 9).Value = 1.038860145888054
$ws.Cells.Item(25, 10).Value = 1.040956739906127
$ws.Cells.Item(25, 11).Value = 1.046764639191818
$ws.Cells.Item(25, 12).Value = 1.037506526197506
$ws.Cells.Item(25, 13).Value = 1.053608226949513
